# GILDNamedTradeLongHold.xlsx - "traded, fixed issues with the repeater"
#
# A new trade was recorded in row 3: the repeater now also stamps the
# trade Date (A), BuyPrice (D) and Holding flag (G) instead of leaving
# them blank (Principle in C3 was already populated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-seed A3/G3 with the same formatting as the analogous cells in row 2
# (copy formats only, first) so the new cells reuse the existing date /
# boolean style rather than minting a fresh cell style.
$ws.Range("A2:A2").Copy($ws.Range("A3"))
$ws.Range("G2:G2").Copy($ws.Range("G3"))

# Now write the actual values for the new trade row.
$ws.Range("A3").Value = [DateTime]::FromOADate(42654.743692129632)
$ws.Range("D3").Value = 75.5
$ws.Range("G3").Value = $true
